$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new entry for 9th March 2025
$ws.Range("B27").Value = 9
$ws.Range("C27").Value = "Reshape the Matrix"
$ws.Range("E27").Value = "LeetCode"

# Update the active selection to match the author's final cursor position
$ws.Range("D32").Select()
